# Shift all Timestamp (col A) values forward by 2 days, and regenerate the
# corresponding "Lookup" (col E) text so it keeps matching the date portion
# of the (now shifted) Timestamp plus the existing Quarter (col D) index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 195

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellD = $ws.Cells.Item($r, 4)
    $cellE = $ws.Cells.Item($r, 5)

    # Shift the timestamp by 2 days (serial date arithmetic).
    $newSerial = $cellA.Value2 + 2
    $cellA.Value2 = $newSerial

    # Recompute the day-count part of the serial (drop the fractional part)
    # and convert back to a real date so we can format it as dd.MM.yyyy.
    $days = [Math]::Floor($newSerial)
    $epoch = Get-Date -Year 1899 -Month 12 -Day 30
    $rowDate = $epoch.AddDays($days)

    $quarter = $cellD.Value2
    $cellE.Value2 = $rowDate.ToString("dd.MM.yyyy") + [string]$quarter
}
